$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column D (Email ID) to hold "Gender"
$ws.Range("D1").EntireColumn.Insert()

# Match the width that was manually set on the new Gender column (close to column C's width)
$ws.Columns(4).ColumnWidth = 9

# Header + values for the new "Gender" column
$ws.Range("D1").Value = "Gender"
$ws.Range("D2").Value = "M"
$ws.Range("D3").Value = "F"
$ws.Range("D4").Value = "F"
$ws.Range("D5").Value = "F"

# The hyperlinks that used to sit on D2:D5 (Email ID column) now need to move to E2:E5.
# Inserting the column doesn't auto-shift the worksheet's hyperlink anchors, so fix them up.
$ws.Range("D2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:jatin@golivefaster.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:divya@golivefaster.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:janhavi@golivefaster.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:shriya@golivefaster.com") | Out-Null

# Re-apply the Hyperlink cell style (Add() bumps these to a fresh style index; restore the
# original shared "Hyperlink" style the cells already had after the column shift).
$ws.Range("E2:E5").Style = "Hyperlink"
